# chore(runtime): publish files + archive (2025-12-18 11:04:48)
# Applies the KHL stats refresh: two new matches on Matches_SOG, the
# as_of_utc rollover (2025-12-16 -> 2025-12-17) across Shots_HA /
# Shots_Summary / Meta_ext, plus the stat deltas produced by that refresh.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Matches_SOG: append the two new games played 2025-12-17.
# ---------------------------------------------------------------------
$wsMatches = $wb.Worksheets.Item("Matches_SOG")

$newMatches = @(
    @{ Row = 381; Uid = "897873"; Date = "2025-12-17T12:15:00"; Home = "Амур";    Away = "Адмирал";     Sog = 36; Soa = 30 },
    @{ Row = 382; Uid = "897874"; Date = "2025-12-17T19:30:00"; Home = "ХК Сочи"; Away = "Металлург Мг"; Sog = 26; Soa = 25 }
)

foreach ($m in $newMatches) {
    $r = $m.Row
    # uid is numeric-looking text in this sheet (see row 2..380) - force text
    # storage via NumberFormat so "897873" isn't coerced into a number.
    $wsMatches.Cells.Item($r, 1).NumberFormat = "@"
    $wsMatches.Cells.Item($r, 1).Value = $m.Uid
    $wsMatches.Cells.Item($r, 2).Value = $m.Date
    $wsMatches.Cells.Item($r, 3).Value = $m.Home
    $wsMatches.Cells.Item($r, 4).Value = $m.Away
    $wsMatches.Cells.Item($r, 5).Value = $m.Sog
    $wsMatches.Cells.Item($r, 6).Value = $m.Soa
    $wsMatches.Cells.Item($r, 7).Value = "khl_text"
}

# ---------------------------------------------------------------------
# 2) Shots_HA: as_of_utc rolls from 2025-12-16T19:30:00Z to
#    2025-12-17T19:30:00Z for every team, plus the per-team home/away
#    shots-on-goal counters that moved because of the two new games.
# ---------------------------------------------------------------------
$wsHA = $wb.Worksheets.Item("Shots_HA")

for ($r = 2; $r -le 23; $r++) {
    $wsHA.Cells.Item($r, 4).Value = "2025-12-17T19:30:00Z"
}

# row 4 = Адмирал (away in game 897873)
$wsHA.Cells.Item(4, 6).Value = 17
$wsHA.Cells.Item(4, 11).Value = 527
$wsHA.Cells.Item(4, 12).Value = 480
$wsHA.Cells.Item(4, 13).Value = 31
$wsHA.Cells.Item(4, 14).Value = 28.2

# row 6 = Амур (home in game 897873)
$wsHA.Cells.Item(6, 5).Value = 19
$wsHA.Cells.Item(6, 7).Value = 584
$wsHA.Cells.Item(6, 8).Value = 648
$wsHA.Cells.Item(6, 9).Value = 30.7
$wsHA.Cells.Item(6, 10).Value = 34.1

# row 13 = Металлург Мг (away in game 897874)
$wsHA.Cells.Item(13, 6).Value = 17
$wsHA.Cells.Item(13, 11).Value = 518
$wsHA.Cells.Item(13, 12).Value = 498
$wsHA.Cells.Item(13, 13).Value = 30.5
$wsHA.Cells.Item(13, 14).Value = 29.3

# row 22 = ХК Сочи (home in game 897874)
$wsHA.Cells.Item(22, 5).Value = 17
$wsHA.Cells.Item(22, 7).Value = 498
$wsHA.Cells.Item(22, 8).Value = 537
$wsHA.Cells.Item(22, 9).Value = 29.3
$wsHA.Cells.Item(22, 10).Value = 31.6

# ---------------------------------------------------------------------
# 3) Shots_Summary: same as_of_utc rollover, plus the combined
#    (home+away) totals for the four teams that played.
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Shots_Summary")

for ($r = 2; $r -le 23; $r++) {
    $wsSummary.Cells.Item($r, 4).Value = "2025-12-17T19:30:00Z"
}

# row 4 = Адмирал
$wsSummary.Cells.Item(4, 5).Value = 33
$wsSummary.Cells.Item(4, 6).Value = 1111
$wsSummary.Cells.Item(4, 7).Value = 913
$wsSummary.Cells.Item(4, 8).Value = 33.7
$wsSummary.Cells.Item(4, 9).Value = 27.7

# row 6 = Амур
$wsSummary.Cells.Item(6, 5).Value = 36
$wsSummary.Cells.Item(6, 6).Value = 1069
$wsSummary.Cells.Item(6, 7).Value = 1292
$wsSummary.Cells.Item(6, 8).Value = 29.7
$wsSummary.Cells.Item(6, 9).Value = 35.9

# row 13 = Металлург Мг
$wsSummary.Cells.Item(13, 5).Value = 34
$wsSummary.Cells.Item(13, 6).Value = 1160
$wsSummary.Cells.Item(13, 7).Value = 928
$wsSummary.Cells.Item(13, 8).Value = 34.1

# row 22 = ХК Сочи
$wsSummary.Cells.Item(22, 5).Value = 33
$wsSummary.Cells.Item(22, 6).Value = 913
$wsSummary.Cells.Item(22, 7).Value = 1125
$wsSummary.Cells.Item(22, 9).Value = 34.1

# ---------------------------------------------------------------------
# 4) Meta_ext: bump as_of_utc + build_version for this publish.
# ---------------------------------------------------------------------
$wsMeta = $wb.Worksheets.Item("Meta_ext")
$wsMeta.Cells.Item(2, 2).Value = "2025-12-17T19:30:00Z"
$wsMeta.Cells.Item(2, 4).Value = 64
